$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: "Judul Dataset" value -> New Plant Diseases Dataset
$ws.Range("A10").Value = "New Plant Diseases Dataset"

# Row 16: "Judul Artikel" value -> new article title
$ws.Range("A16").Value = "KLASIFIKASI PENYAKIT TANAMAN APEL DARI CITRA DAUN DENGAN CONVOLUTIONAL NEURAL NETWORK"

# Row 19: "Link Artikel" value -> new article link, turned into a real hyperlink
$ws.Range("A19").Value = "https://jurnal.wicida.ac.id/index.php/sebatik/article/download/1060/297/"
$ws.Hyperlinks.Add($ws.Range("A19"), "https://jurnal.wicida.ac.id/index.php/sebatik/article/download/1060/297/")

# Row 13: "Link Dataset" value -> new dataset link, styled like a hyperlink but not an
# actual clickable hyperlink (reuses the Hyperlink style minted above).
$ws.Range("A13").Value = "https://www.kaggle.com/vipoooool/new-plant-diseases-dataset"
$ws.Range("A13").Style = "Hyperlink"

# Move the active selection to A10, matching the saved cursor position.
$ws.Range("A10").Select() | Out-Null
